# Add a new "2021" column (S) to the right of the existing "2020" column (R),
# copying the source column's formatting cell-by-cell and then writing the
# new values, then update the active selection to match the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (R source cell, new value for S)
$rows = @(
    @{ Row = 4;  Value = 2021 },
    @{ Row = 5;  Value = 6.1 },
    @{ Row = 6;  Value = 1.6 },
    @{ Row = 7;  Value = 3.6 },
    @{ Row = 8;  Value = 27.2 },
    @{ Row = 9;  Value = 7.2 },
    @{ Row = 10; Value = 2.6 },
    @{ Row = 11; Value = 12.5 },
    @{ Row = 12; Value = 6.4 },
    @{ Row = 13; Value = 5.2 },
    @{ Row = 14; Value = 0.9 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $src = $ws.Range("R$r")
    $dst = $ws.Range("S$r")
    $src.Copy($dst)
    $dst.Value = $item.Value
}

# Restore the active selection described by the edit.
$ws.Range("Q19").Select()
